$d = $word.ActiveDocument

# Locate the paragraph that holds "Kem cho maaz ama" (it currently ends
# with a trailing space run) and the paragraph right after it, which
# holds "Chalo chalo".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Kem cho maaz ama") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Remove the whole following paragraph ("Chalo chalo"), including its
    # paragraph mark, merging it away and leaving the previous paragraph's
    # own mark/properties intact.
    $next = $d.Paragraphs.Item($target + 1)
    $next.Range.Delete()

    # Trim the trailing space left at the end of the "Kem cho maaz ama"
    # paragraph. Paragraph.Range.Text includes the trailing paragraph
    # mark (Chr(13)), so the character to inspect is the one just before
    # that mark.
    $p = $d.Paragraphs.Item($target)
    $r = $p.Range
    $textLen = $r.Text.Length
    if ($textLen -gt 1 -and $r.Text.Substring($textLen - 2, 1) -eq " ") {
        $trim = $d.Range($r.Start + $textLen - 2, $r.Start + $textLen - 1)
        $trim.Delete()
    }
}
